# Updates cryptos list values (prices & volume %) and restores the
# HuobiToken / MXToken / ARBITRUM row order to match the refreshed ranking.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.072.63'
$ws.Range('E2').Value = '  -0.44%  '
$ws.Range('D3').Value = '1.652.21'
$ws.Range('E3').Value = '  -0.38%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').Value = "'217.51"
$ws.Range('E5').Value = '  +0.47%  '
$ws.Range('D6').Value = "'0.5251"
$ws.Range('E6').Value = '  +2.29%  '
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('D8').Value = "'0.2592"
$ws.Range('E8').Value = '  -1.65%  '
$ws.Range('D9').Value = "'0.06337"
$ws.Range('E9').Value = '  +1.19%  '
$ws.Range('D10').Value = "'20.34"
$ws.Range('E10').Value = '  -2.08%  '
$ws.Range('D11').Value = "'0.07801"
$ws.Range('E11').Value = '  +0.64%  '
$ws.Range('D12').Value = "'4.494"
$ws.Range('E12').Value = '  +1.00%  '
$ws.Range('D13').Value = '1.662.13'
$ws.Range('E13').Value = '  +0.15%  '
$ws.Range('E14').Value = '  +1.02%  '
$ws.Range('D15').Value = '0.0₅8228'
$ws.Range('E15').Value = '  +1.81%  '
$ws.Range('D16').Value = "'65.38"
$ws.Range('E16').Value = '  +0.85%  '
$ws.Range('D17').Value = '26.101.12'
$ws.Range('E17').Value = '  -0.34%  '
$ws.Range('D18').Value = "'1.003"
$ws.Range('E18').Value = '  -0.10%  '
$ws.Range('D19').Value = "'4.571"
$ws.Range('E19').Value = '  -0.93%  '
$ws.Range('D20').Value = "'190.77"
$ws.Range('E20').Value = '  -0.65%  '
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('D22').Value = "'6.029"
$ws.Range('E22').Value = '  +0.59%  '
$ws.Range('D23').Value = "'1.004"
$ws.Range('E23').Value = '  -0.13%  '
$ws.Range('D24').Value = "'142.21"
$ws.Range('E24').Value = '  +1.54%  '
$ws.Range('D25').Value = "'0.1231"
$ws.Range('E25').Value = '  +1.09%  '
$ws.Range('D26').Value = "'7.237"
$ws.Range('E26').Value = '  +0.24%  '
$ws.Range('E27').Value = '  -0.49%  '
$ws.Range('E28').Value = '  -0.35%  '
$ws.Range('D29').Value = "'0.05834"
$ws.Range('E29').Value = '  -1.62%  '
$ws.Range('E30').Value = '  +0.15%  '
$ws.Range('D31').Value = "'3.544"
$ws.Range('E31').Value = '  -0.33%  '
$ws.Range('D32').Value = "'3.254"
$ws.Range('E32').Value = '  +0.15%  '
$ws.Range('D33').Value = "'1.583"
$ws.Range('E33').Value = '  -0.70%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').Value = "'2.413"
$ws.Range('E34').Value = '  -0.42%  '
$ws.Range('B35').Value = 'MXToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D35').Value = "'2.779"
$ws.Range('E35').Value = '  +0.31%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').Value = "'0.9447"
$ws.Range('E36').Value = '  -1.90%  '
$ws.Range('D37').Value = "'0.5725"
$ws.Range('E37').Value = '  +1.54%  '
$ws.Range('D38').Value = "'0.01609"
$ws.Range('E38').Value = '  +1.17%  '
$ws.Range('D39').Value = "'5.754"
$ws.Range('E39').Value = '  -3.27%  '
$ws.Range('D40').Value = "'0.8438"
$ws.Range('E40').Value = '  -1.80%  '
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('D42').Value = "'103.45"
$ws.Range('E42').Value = '  +3.14%  '
$ws.Range('D43').Value = '1.025.16'
$ws.Range('E43').Value = '  +1.64%  '
$ws.Range('D44').Value = '1.796.89'
$ws.Range('E45').Value = '  +0.76%  '
$ws.Range('D46').Value = "'0.9999"
$ws.Range('E46').Value = '  -0.49%  '
$ws.Range('D47').Value = "'0.4315"
$ws.Range('E47').Value = '  +3.14%  '
$ws.Range('D48').Value = "'0.05147"
$ws.Range('E48').Value = '  -0.38%  '
$ws.Range('D49').Value = "'1.463"
$ws.Range('E49').Value = '  +1.15%  '
$ws.Range('E50').Value = '  -2.55%  '
$ws.Range('D51').Value = "'0.09638"
